$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 764204132003
$ws.Range("B2").Value = "In The Zone - 3 Items"
$ws.Range("C2").Value = 175540207634
$ws.Range("D2").Value = 1599416074258
$ws.Range("E2").Value = "In The Zone - 3 Items"

$ws.Range("A3").Value = 764204132010
$ws.Range("B3").Value = "In The Zone - 5 Items"
$ws.Range("C3").Value = 175535685650
$ws.Range("D3").Value = 1599409848338
$ws.Range("E3").Value = "In The Zone - 5 Items"

$ws.Range("A3").Select()
